$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 14-17 (old sheet had 17 data rows, new sheet only needs 13)
$ws.Range("A14:F17").EntireRow.Delete() | Out-Null

# Clear the whole data range first so stale values (e.g. leftover F column
# messages) don't linger, then refill with the new values.
$ws.Range("A2:F13").ClearContents() | Out-Null

$empId = "1126673"

$data = @(
    @{ B = "test_factorial_num_case1";   C = "factorial_num";    D = "passed"; E = 0.00038909912109375;  F = $null },
    @{ B = "test_factorial_num_case2";   C = "factorial_num";    D = "passed"; E = 0.0003588199615478516; F = $null },
    @{ B = "test_factorial_num_case3";   C = "factorial_num";    D = "passed"; E = 0.0003521442413330078; F = $null },
    @{ B = "test_factorial_num_case4";   C = "factorial_num";    D = "passed"; E = 0.0003654956817626953; F = $null },
    @{ B = "test_factorial_num_case5";   C = "factorial_num";    D = "failed"; E = 0.0004706382751464844; F = "TypeError: reduce() of empty sequence with no initial value" },
    @{ B = "test_factorial_num_case6";   C = "factorial_num";    D = "failed"; E = 0.0003962516784667969; F = "Failed: DID NOT RAISE <class 'ValueError'>" },
    @{ B = "test_even_pos_upcase_case1"; C = "even_pos_upcase";  D = "passed"; E = 0.0003616809844970703; F = $null },
    @{ B = "test_even_pos_upcase_case2"; C = "even_pos_upcase";  D = "passed"; E = 0.0003647804260253906; F = $null },
    @{ B = "test_even_pos_upcase_case3"; C = "even_pos_upcase";  D = "passed"; E = 0.0004978179931640625; F = $null },
    @{ B = "test_even_pos_upcase_case4"; C = "even_pos_upcase";  D = "passed"; E = 0.0003888607025146485; F = $null },
    @{ B = "test_even_pos_upcase_case5"; C = "even_pos_upcase";  D = "passed"; E = 0.0004103183746337891; F = $null },
    @{ B = "test_even_pos_upcase_case6"; C = "even_pos_upcase";  D = "passed"; E = 0.0003874301910400391; F = $null }
)

$row = 2
foreach ($rec in $data) {
    $aCell = $ws.Cells.Item($row, 1)
    $aCell.NumberFormat = "@"
    $aCell.Value = $empId
    $ws.Cells.Item($row, 2).Value = $rec.B
    $ws.Cells.Item($row, 3).Value = $rec.C
    $ws.Cells.Item($row, 4).Value = $rec.D
    $ws.Cells.Item($row, 5).Value = $rec.E
    if ($rec.F) {
        $ws.Cells.Item($row, 6).Value = $rec.F
    }
    $row++
}
